$d = $word.ActiveDocument

$ids = @("p078v_4", "p079r_1", "p079r_2", "p079r_3")

foreach ($id in $ids) {
    $old = "<id>" + $id + "</id>"
    $new = "<id>" + $id + "</id>"
    $found = $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
    Write-Host "Replaced" $id ":" $found
}
